# Weekly update: a new daily price record was inserted for
# "Mercado Mayorista Lo Valledor de Santiago - Frambuesa" ahead of the
# existing row 100, pushing all subsequent rows (old 100..166) down by
# one position (new 101..167).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 100 - Excel shifts every row
# at/after 100 down by one (old row100 -> row101, ..., old row166 -> row167)
# and the sheet's used range grows from T166 to T167 automatically.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record's data.
$ws.Cells.Item(100, 1).Value = 6
$ws.Cells.Item(100, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(100, 3).Value = "Metropolitana"
$ws.Cells.Item(100, 4).Value = 44603
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100101
$ws.Cells.Item(100, 8).Value = "Berries"
$ws.Cells.Item(100, 9).Value = 100101004
$ws.Cells.Item(100, 10).Value = "Frambuesa"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 230
$ws.Cells.Item(100, 14).Value = 6000
$ws.Cells.Item(100, 15).Value = 6000
$ws.Cells.Item(100, 16).Value = 6000
$ws.Cells.Item(100, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(100, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(100, 19).Value = 3000
$ws.Cells.Item(100, 20).Value = 2
